$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "SXT"
$ws.Range("B11").Value = 8653
$ws.Range("C11").Value = 11086.87949899516
$ws.Range("D11").Value = 0.6457036000618893
